# Aging stock summary: re-bucket the aging days, swap several brand/item
# rows around, and drop the old last data row (Toti moves up into the
# "61 - 90 Days" bucket instead of having its own trailing row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: first row just changes its bucket label to "Expired" ---
$ws.Range("A2").Value = "Expired"

# --- Row 3: was Naprox / 1-15 Days -> becomes Flucloxin / 16-30 Days ---
$ws.Range("A3").Value = "16 - 30 Days"
$ws.Range("B3").Value = "Flucloxin"
$ws.Range("D3").Value = "Flucloxin 100ml Dry Suspension"
$ws.Range("Y3").ClearContents()
$ws.Range("G3").Value = 22

# --- Row 4: was Toperin / 1-15 Days -> becomes Flucloxin / 16-30 Days ---
$ws.Range("A4").Value = "16 - 30 Days"
$ws.Range("B4").Value = "Flucloxin"
$ws.Range("D4").Value = "Flucloxin 500mg Capsule 40's"
$ws.Range("H4").ClearContents()
$ws.Range("G4").Value = 7
$ws.Range("Y4").Value = 1

# --- Row 5: was Zithrox 50ml / 1-15 Days -> becomes Mebidal / 16-30 Days ---
$ws.Range("A5").Value = "16 - 30 Days"
$ws.Range("B5").Value = "Mebidal"
$ws.Range("D5").Value = "Mebidal Tablet"
$ws.Range("U5").ClearContents()
$ws.Range("AC5").Value = 24

# --- Row 6: was Zithrox 20ml / 1-15 Days -> becomes Naprox / 16-30 Days ---
$ws.Range("A6").Value = "16 - 30 Days"
$ws.Range("B6").Value = "Naprox"
$ws.Range("D6").Value = "Naprox Plus 500mg Tablet - 36's"
$ws.Range("M6").ClearContents()
$ws.Range("Y6").Value = 1

# --- Row 7: was Flucloxin 100ml / 31-60 Days -> becomes Osticare / 16-30 Days ---
$ws.Range("A7").Value = "16 - 30 Days"
$ws.Range("B7").Value = "Osticare"
$ws.Range("D7").Value = "Osticare Tablet 30's"
$ws.Range("G7").ClearContents()
$ws.Range("J7").Value = 1
$ws.Range("Y7").Value = 1

# --- Row 8: was Flucloxin 500mg / 31-60 Days (bucket stays) -> becomes Oradin ---
$ws.Range("B8").Value = "Oradin"
$ws.Range("D8").Value = "Oradin 60ml Suspension"
$ws.Range("G8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("J8").Value = 2

# --- Row 9: was Mebidal / 31-60 Days -> becomes Dinafex / 61-90 Days ---
$ws.Range("A9").Value = "61 - 90 Days"
$ws.Range("B9").Value = "Dinafex"
$ws.Range("D9").Value = "Dinafex 50ml Suspension"
$ws.Range("AC9").Value = 1

# --- Row 10: was Oradin / 31-60 Days -> becomes 61-90 Days (brand/item same) ---
$ws.Range("A10").Value = "61 - 90 Days"
$ws.Range("J10").Value = 1

# --- Row 11: was Osticare / 31-60 Days -> becomes Toti / 61-90 Days ---
$ws.Range("A11").Value = "61 - 90 Days"
$ws.Range("B11").Value = "Toti"
$ws.Range("D11").Value = "Toti 100ml Syrup"
$ws.Range("Y11").ClearContents()
$ws.Range("G11").Value = 14
$ws.Range("J11").Value = 5
$ws.Range("AA11").Value = 33

# --- Row 12: was Oradin / 61-90 Days -> becomes Zithrox (bucket stays) ---
$ws.Range("B12").Value = "Zithrox"
$ws.Range("D12").Value = "Zithrox 35ml Dry Suspension"
$ws.Range("J12").Value = 7
$ws.Range("Y12").Value = 1

# --- Old row 13 (Toti / 61-90 Days) is folded into row 11 above; remove it ---
$ws.Rows.Item(13).Delete()
